# Added two new Mac-Addresses: append 10 new device rows (3000166-3000175)
# for regcntr_id 10001 to the master-reg_center_device_h sheet, following
# the same pattern as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startDeviceId = 3000166
$rowCount = 10

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $deviceId = $startDeviceId + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Reflect the post-edit view state: selection on the last-entered cell,
# scrolled so the new rows are visible.
$win = $excel.ActiveWindow
$win.ScrollRow = 140
$ws.Range("E155").Select() | Out-Null
